$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: blank separator row - bold + text number format, no fill (creates a new style)
foreach ($col in @("A","B","C","E")) {
    $cell = $ws.Range($col + "8")
    $cell.Font.Bold = $true
    $cell.NumberFormat = "@"
}

# Row 10: header row for new "Acura" scenario (reuse row 4's formatting)
$ws.Range("A4:E4").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A10").Value = "Acura"
$ws.Range("B10").Value = "stockType"
$ws.Range("C10").Value = "make"
$ws.Range("D10").Value = "expectedModels"
$ws.Range("E10").Value = "Acura"

# Row 11: data row for new "Acura" scenario (reuse row 5's formatting)
$ws.Range("B5:D5").Copy()
$ws.Range("B11:D11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B11").Value = "New"
$ws.Range("C11").Value = "Acura"
$ws.Range("D11").Value = "All models,ILX,Integra,MDX,NSX,RDX,TLX"

$excel.CutCopyMode = $false

# Update selection to mirror what Excel leaves as "active cell" after editing
$ws.Range("D12").Select() | Out-Null
